$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.553.28'
$ws.Range('E2').Value = '  -1.46%  '

$ws.Range('D3').Value = '2.900.06'
$ws.Range('E3').Value = '  -2.43%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.42'
$ws.Range('E5').Value = '  -2.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.76'
$ws.Range('E6').Value = '  -4.35%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.548'
$ws.Range('E8').Value = '  -3.26%  '

$ws.Range('D9').Value = '2.908.07'
$ws.Range('E9').Value = '  -2.47%  '

$ws.Range('E10').Value = '  -5.17%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.97'
$ws.Range('E11').Value = '  -2.54%  '

$ws.Range('E12').Value = '  -2.62%  '

$ws.Range('D13').Value = '3.409.52'
$ws.Range('E13').Value = '  -2.28%  '

$ws.Range('E14').Value = '  +2.36%  '

$ws.Range('D15').Value = '60.552.60'
$ws.Range('E15').Value = '  -1.65%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.55'
$ws.Range('E16').Value = '  -4.38%  '

$ws.Range('D17').Value = '2.906.44'
$ws.Range('E17').Value = '  -2.30%  '

$ws.Range('E18').Value = '  -3.83%  '

$ws.Range('E19').Value = '  -3.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.60'
$ws.Range('E20').Value = '  -3.64%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.36'
$ws.Range('E21').Value = '  -7.28%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.53'
$ws.Range('E22').Value = '  -2.30%  '

$ws.Range('E23').Value = '  -0.24%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.71'
$ws.Range('E24').Value = '  +0.97%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.61'
$ws.Range('E25').Value = '  -1.54%  '

$ws.Range('E26').Value = '  -3.75%  '

$ws.Range('E27').Value = '  -5.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.80'
$ws.Range('E29').Value = '  -4.78%  '

$ws.Range('D30').Value = '0.0₃0849'
$ws.Range('E30').Value = '  -9.24%  '

$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('E32').Value = '  -2.51%  '

$ws.Range('E33').Value = '  -3.83%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.23'
$ws.Range('E34').Value = '  -5.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.31'
$ws.Range('E35').Value = '  -6.90%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.56'
$ws.Range('E36').Value = '  -5.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.997'
$ws.Range('E37').Value = '  -6.74%  '

$ws.Range('E38').Value = '  -5.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.66'
$ws.Range('E39').Value = '  +0.43%  '

$ws.Range('E40').Value = '  -5.04%  '

$ws.Range('E41').Value = '  -5.28%  '

$ws.Range('D42').Value = '2.290.17'
$ws.Range('E42').Value = '  -5.02%  '

$ws.Range('E43').Value = '  -3.26%  '

$ws.Range('E44').Value = '  -1.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.41'
$ws.Range('E45').Value = '  -7.56%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.94'
$ws.Range('E47').Value = '  -2.95%  '

$ws.Range('E48').Value = '  -3.28%  '

$ws.Range('E49').Value = '  -1.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0919'
$ws.Range('E50').Value = '  -3.30%  '

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.27'
$ws.Range('E51').Value = '  -7.32%  '

